$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the existing "dsVal1" column (F), copying
# formatting from the column to the LEFT (xlFormatFromLeftOrAbove = 0) so the
# new header cells (row 1) pick up the same header style as the rest of row 1,
# and the new row-2 cells stay unstyled like their D2/E2 neighbours.
$ws.Columns("F:G").Insert(0)

# New header labels for the inserted columns.
$ws.Range("F1").Value = "xrayTestExecKey"
$ws.Range("G1").Value = "xrayTestKey"

# New data values for row 2 under the inserted columns.
$ws.Range("F2").Value = "n/a"
$ws.Range("G2").Value = "n/a"

# The cell-level style objects created by Insert() for the new row-2 cells
# differ subtly from an untouched cell; line them back up with the plain,
# unstyled neighbour cell (D2) so they match the rest of that row.
$ws.Range("F2:G2").Style = $ws.Range("D2").Style

# The old, empty dsVal1 cell (row 2) shifted from F2 to H2 when the columns
# were inserted; it never held any content, so drop the now-stray empty cell
# entirely instead of leaving a placeholder behind.
$ws.Range("H2").ClearContents()

# Row-height tweaks that came with the new layout.
$ws.Rows("1").RowHeight = 28.2
$ws.Rows("2").RowHeight = 15

# The shifted dsVal1/dsVal2 columns (now H:I) end up a touch narrower than
# the newly inserted columns, and dsVal3 (now J) narrower still - restore
# that width grouping instead of leaving all of E:H at the copied width.
$ws.Columns("H:I").ColumnWidth = 15
$ws.Columns("J").ColumnWidth = 14.5

# Selection ends up on H2 (first cell of the original dsVal1 column, now
# shifted right) after the edit.
$ws.Range("H2").Select() | Out-Null
